$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Anchor" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Anchor"

# Populate the data (header + rows) as literal text values
$data = @(
    @("Anchor Text", "Count", "Percentage", "Link Type"),
    @("Hosn Insurance", "234", "18.2%", "Branded"),
    @("car insurance UAE", "98", "7.6%", "Exact Match"),
    @("click here", "156", "12.1%", "Generic"),
    @("https://hosninsurance.ae", "89", "6.9%", "URL"),
    @("best insurance company", "67", "5.2%", "Partial Match"),
    @("Dubai insurance", "54", "4.2%", "Partial Match"),
    @("auto insurance", "43", "3.3%", "Exact Match"),
    @("visit website", "87", "6.8%", "Generic"),
    @("Hosn Insurance Dubai", "76", "5.9%", "Branded"),
    @("insurance services", "45", "3.5%", "Partial Match"),
    @("read more", "123", "9.6%", "Generic"),
    @("vehicle insurance UAE", "32", "2.5%", "Exact Match"),
    @("UAE insurance provider", "28", "2.2%", "Partial Match"),
    @("learn more", "91", "7.1%", "Generic"),
    @("comprehensive coverage", "21", "1.6%", "Partial Match"),
    @("motor insurance", "19", "1.5%", "Exact Match"),
    @("insurance quotes", "17", "1.3%", "Partial Match"),
    @("get quote", "25", "1.9%", "Generic"),
    @("Hosn", "35", "2.7%", "Branded"),
    @("homepage", "14", "1.1%", "Generic"),
    @("insurance brokers", "12", "0.9%", "Partial Match"),
    @("contact us", "18", "1.4%", "Generic"),
    @("about company", "16", "1.2%", "Generic"),
    @("UAE car insurance", "13", "1.0%", "Exact Match"),
    @("professional services", "11", "0.9%", "Partial Match"),
)

# Columns B (Count) and C (Percentage) hold digit-only / percent-looking
# text that Excel would otherwise auto-convert to numbers, so force those
# two columns to Text format before writing so the values stay literal
# strings (matching the source report's inline-string cells).
$ws.Range("B1:C26").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $cell = $ws.Cells.Item($i + 1, $j + 1)
        $cell.Value = $row[$j]
    }
}
